# feat!: removal of option `fieldMatchType`
# new default is `labelTypeBrackets` to avoid any collisions with labels and types
#
# The header row now encodes both the label and the matched field type using
# the `Label[type]` bracket notation instead of a plain label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID[product_ID]"
$ws.Range("B1").Value = "Quantity[quantity]"
$ws.Range("C1").Value = "ProductTitle[title]"
$ws.Range("D1").Value = "UnitPrice[price]"
$ws.Range("E1").Value = "validFrom[validFrom]"
$ws.Range("F1").Value = "timestamp[timestamp]"
$ws.Range("G1").Value = "date[date]"
$ws.Range("H1").Value = "time[time]"
$ws.Range("I1").Value = "WRONGCOLUMN[TEST]"

$ws.Range("A1:H1").Select()
$ws.Range("H1").Activate()
